# Generate Report for Archive
#
# The localization run moved this document out of the "Ready for handoff"
# state into "In Translation": update the Status value everywhere it is
# shown (the Overview rollup columns for zh-cn/de-de, plus each locale
# sheet's own Status column) and re-fit the columns that used to hold the
# old, longer status text so the refreshed report renders tidily.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: one rollup column per locale (E = zh-cn, F = de-de) ---
if ($overview.Range("E2").Value() -eq $oldStatus) { $overview.Range("E2").Value = $newStatus }
if ($overview.Range("F2").Value() -eq $oldStatus) { $overview.Range("F2").Value = $newStatus }

# --- Per-locale sheets: Status is column C ---
if ($zhcn.Range("C2").Value() -eq $oldStatus) { $zhcn.Range("C2").Value = $newStatus }
if ($dede.Range("C2").Value() -eq $oldStatus) { $dede.Range("C2").Value = $newStatus }

# Re-fit the columns that used to hold the longer "Ready for handoff" text
# now that the shorter "In Translation" status is in place.
$newWidth = 12.5
$overview.Range("E:E").ColumnWidth = $newWidth
$overview.Range("F:F").ColumnWidth = $newWidth
$zhcn.Range("C:C").ColumnWidth = $newWidth
$dede.Range("C:C").ColumnWidth = $newWidth
